$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers -----------------------------------------------------
# Existing columns A-E (date, cpu, netin, netout, scale) remain unchanged.
# Column F stays empty, as before.

$ws.Range("G1").Value = "up_cpu"
$ws.Range("H1").Value = "down_cpu"
$ws.Range("I1").Value = "up_netin"
$ws.Range("J1").Value = "down_netin"
$ws.Range("K1").Value = "up_netout"
$ws.Range("L1").Value = "down_netout"
$ws.Range("M1").Value = "cpu_acc"
$ws.Range("N1").Value = "cpu_model"
$ws.Range("O1").Value = "cpu_pred1"
$ws.Range("P1").Value = "cpu_pred2"
# Q1 held "netin_pred2" before this edit; it is not reused by any new
# header, so clear it out.
$ws.Range("Q1").ClearContents()
$ws.Range("R1").Value = "cpu_pred3"
$ws.Range("S1").Value = "netin_acc"
$ws.Range("T1").Value = "netin_model"
$ws.Range("U1").Value = "netin_pred1"
$ws.Range("V1").Value = "netin_pred2"
$ws.Range("W1").Value = "netin_pred3"
$ws.Range("X1").Value = "netout_acc"
$ws.Range("Y1").Value = "netout_model"
$ws.Range("Z1").Value = "netout_pred1"
$ws.Range("AA1").Value = "netout_pred2"
$ws.Range("BB1").Value = "netout_pred3"

# --- Row 2 data ----------------------------------------------------------
# Plain numeric cells (stay numeric)
$ws.Range("A2").Value = "05/08/2021 02:45:17"
$ws.Range("B2").Value = 0.3390000000000001
$ws.Range("C2").Value = 14.368
$ws.Range("D2").Value = 12.736
$ws.Range("E2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("V2").Value = 0

# Cells that hold numeric-looking text - force Text format so the value
# is written as a string, not re-parsed as a number.
$textCells = "G2","H2","I2","J2","K2","L2","N2","O2","P2","S2","T2","U2","X2","Y2","Z2"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("G2").Value = "0.3390000000000001"
$ws.Range("H2").Value = "0.3390000000000001"
$ws.Range("I2").Value = "14.368"
$ws.Range("J2").Value = "14.368"
$ws.Range("K2").Value = "12.736"
$ws.Range("L2").Value = "12.736"
$ws.Range("M2").Value = "(0, 0, 0)"
$ws.Range("N2").Value = "0"
$ws.Range("O2").Value = "0"
$ws.Range("P2").Value = "0"
$ws.Range("R2").Value = "(0, 0, 0)"
$ws.Range("S2").Value = "0"
$ws.Range("T2").Value = "0"
$ws.Range("U2").Value = "0"
$ws.Range("W2").Value = "(0, 0, 0)"
$ws.Range("X2").Value = "0"
$ws.Range("Y2").Value = "0"
$ws.Range("Z2").Value = "0"
